$wb = $excel.ActiveWorkbook

# --- Update the narrative text on "Hoja1" (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$a1 = $wsHoja1.Range("A1")
$text = $a1.Value()
$text = $text.Replace("1000 Bs = 4.55 = 17954.55 pesos", "1000 Bs = 4.52 = 17812.22 pesos")
$text = $text.Replace("17954.55 pesos = 4.52 = 958.3 Bs", "17812.22 pesos = 4.5 = 967.81 Bs")
$a1.Value = $text

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 221
$wsTasas.Range("O10").Value = 3936.5
$wsTasas.Range("N12").Value = 3957
$wsTasas.Range("O12").Value = 215
